$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing three rows (former gamma1/2nd-strain, nu duplicates) -
# the refreshed "toybox" calibration only has six fitted parameters.
$ws.Rows("8:10").Delete()

# Row 2: Lambda -> gamma1
$ws.Range("B2").Value = "gamma1"
$ws.Range("C2").Value = 0.01
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "( 0 , 0.02 )"
$ws.Range("F2").Value = 2624.7
$ws.Range("G2").Value = 1

# Row 3: beta1 -> gamma2
$ws.Range("B3").Value = "gamma2"
$ws.Range("C3").Value = 0.09
$ws.Range("D3").Value = 0.09
$ws.Range("E3").Value = "( 0.05 , 0.1 )"
$ws.Range("F3").Value = 2923.43
$ws.Range("G3").Value = 1

# Row 4: beta2 -> nu
$ws.Range("B4").Value = "nu"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "( 0 , 0 )"
$ws.Range("F4").Value = 1693.7
$ws.Range("G4").Value = 1

# Row 5: mu -> rho
$ws.Range("B5").Value = "rho"
$ws.Range("C5").Value = 0.33
$ws.Range("D5").Value = 0.32
$ws.Range("E5").Value = "( 0.06 , 0.7 )"
$ws.Range("F5").Value = 2425.17
$ws.Range("G5").Value = 1

# Row 6: gamma1 -> phi1
$ws.Range("B6").Value = "phi1"
$ws.Range("C6").Value = 0.25
$ws.Range("D6").Value = 0.25
$ws.Range("E6").Value = "( 0.19 , 0.33 )"
$ws.Range("F6").Value = 2516.06
$ws.Range("G6").Value = 1

# Row 7: gamma2 -> phi2
$ws.Range("B7").Value = "phi2"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = "( 0 , 0 )"
$ws.Range("F7").Value = 3197.5
$ws.Range("G7").Value = 1
